$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("A2").Value = 8
$ws.Range("C2").Value = "UTS axial"

# Add a new row 3, mirroring row 2's structure but with updated id/property
$ws.Range("A3").Value = 12
$ws.Range("B3").Value = "A08"
$ws.Range("C3").Value = "UTS axial"
$ws.Range("D3").Value = "Mechanical Properties"
$ws.Range("E3").Value = "220 IPHWR"
$ws.Range("F3").Value = "RAPS-1"
$ws.Range("G3").Value = "'2023"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = "'100"
$ws.Range("H3").Style = "Normal"
$ws.Range("I3").Value = "'10.5"
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = "a"
$ws.Range("K3").Value = "'07-10-2024"
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").Value = "Test entry"

$ws.Range("M3").Value = "Test Value 1"
$ws.Range("N3").Value = "Test Value 2"
$ws.Range("O3").Value = "Test Value 3"
$ws.Range("P3").Value = "Test Value 4"
$ws.Range("Q3").Value = "Test Value 5"
$ws.Range("R3").Value = "Test Value 6"
$ws.Range("S3").Value = "Test Value 7"
$ws.Range("T3").Value = "Test Value 8"
$ws.Range("U3").Value = "Test Value 9"
$ws.Range("V3").Value = "Test Value 10"
$ws.Range("W3").Value = "Test Value 11"
$ws.Range("X3").Value = "Test Value 12"
$ws.Range("Y3").Value = "Test Value 13"
$ws.Range("Z3").Value = "Test Value 14"
$ws.Range("AA3").Value = "Test Value 15"
$ws.Range("AB3").Value = "Test Value 16"
$ws.Range("AC3").Value = "Test Value 17"
$ws.Range("AD3").Value = "Test Value 18"
$ws.Range("AE3").Value = "Test Value 19"
$ws.Range("AF3").Value = "Test Value 20"
$ws.Range("AG3").Value = "Test Value 21"
$ws.Range("AH3").Value = "Test Value 22"
$ws.Range("AI3").Value = "Test Value 23"
$ws.Range("AJ3").Value = "Test Value 24"
